# Update gh-pages to output generated at 456a3b4
# Applies updated "想去人数" (interest count) values in column F
# on the "展览" and "全部类型" worksheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 6715
$ws1.Range("F5").Value = 68
$ws1.Range("F15").Value = 1454
$ws1.Range("F17").Value = 3366
$ws1.Range("F19").Value = 224
$ws1.Range("F21").Value = 2004
$ws1.Range("F22").Value = 115

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6716
$ws4.Range("F5").Value = 68
$ws4.Range("F16").Value = 1454
$ws4.Range("F18").Value = 3366
$ws4.Range("F20").Value = 224
$ws4.Range("F22").Value = 2004
$ws4.Range("F23").Value = 115
